$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

# Row 2: quantity -> 0, serial list cleared
$ws.Range("B2").Value = 0
$ws.Range("C2").ClearContents()

# Row 3: indoor model text updated, quantity -> 0, serial cleared
$ws.Range("D3").Value = "asdf"
$ws.Range("E3").Value = 0
$ws.Range("F3").ClearContents()
